$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel's automatic
# "looks like a date" conversion (e.g. "2025-03-30" -> date serial).
# We write it as a string-literal formula (="...") so it is parsed as text,
# then Copy + PasteSpecial(xlPasteValues=-4163) to collapse the formula down
# to a plain literal value/text cell (no lingering formula, no style churn).
function Set-LiteralText($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Row 2
$ws.Range("A2").Value = "Blah blah blah"
$ws.Range("B2").Value = 4000
Set-LiteralText $ws.Range("C2") "2025-03-30"

# Row 3
$ws.Range("A3").Value = "entertainment"
$ws.Range("B3").Value = 2000
Set-LiteralText $ws.Range("C3") "2025-03-30"

# Row 4
$ws.Range("A4").Value = "cake"
$ws.Range("B4").Value = 200
Set-LiteralText $ws.Range("C4") "2025-03-30"

# Row 5
$ws.Range("A5").Value = "Cards"
$ws.Range("B5").Value = 20000
Set-LiteralText $ws.Range("C5") "2024-03-30"

# Row 6
$ws.Range("A6").Value = "music"
$ws.Range("B6").Value = 20000
Set-LiteralText $ws.Range("C6") "2024-03-30"

# Row 7 (highlighted row - yellow fill)
$ws.Range("A7").Value = "munchee nice"
$ws.Range("B7").Value = 20000
Set-LiteralText $ws.Range("C7") "2024-03-30"
$ws.Range("A7:C7").Interior.Color = 65535

# Row 8 (originally row 2's data: Food / 1500 / 2024-03-29, now pushed down)
$ws.Range("A8").Value = "Food"
$ws.Range("B8").Value = 1500
Set-LiteralText $ws.Range("C8") "2024-03-29"

# Column widths (approximate AutoFit of the final data)
$ws.Columns("A").ColumnWidth = 13.19921875
$ws.Columns("B").ColumnWidth = 7.69921875
$ws.Columns("C").ColumnWidth = 10.09765625

# Selection / view state
[void]$ws.Range("C5").Select()
